# Ajustando Tela de Login
#
# The paragraph describing the community features was split across three
# separate runs (likely from incremental edits in Word). Collapse it back
# into a single run with a single <w:t> by replacing the full sentence
# with itself - Word's Find/Replace rebuilds the matched range as one run
# when it applies the replacement.

$d = $word.ActiveDocument

$old = "A comunidade incluirá recursos como fóruns de discussão, resenhas e avaliações de animes e mangás, páginas de resumo para programas populares, e postagens de usuários que compartilham conteúdo relacionado a animes e mangás. Também incluirá atualizações semanais sobre os lançamentos mais recentes além de informações sobre eventos e convenções de animes e mangás."

$rng = $d.Content
$rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $old, 2) | Out-Null
